$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-06 16:11:42"

$wsZhCn.Range("H2").Value = "2016-09-06 16:11:34"
$wsZhCn.Range("K2").Value = "2016-09-06 16:12:47"

$wsDeDe.Range("H2").Value = "2016-09-06 16:11:42"
$wsDeDe.Range("K2").Value = "2016-09-06 16:12:56"
